$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "Tarun"
$ws.Range("B4").Value = "Nice1234"

$ws.Range("A5").Value = "Tarun"
$ws.Range("B5").Value = "test12"

$ws.Range("A6").Value = "Tarun"
$ws.Range("B6").Value = "test12"

$ws.Range("A7").Value = "Nice"
$ws.Range("B7").Value = "Nice"

$ws.Range("A8").Value = "Nice"
$ws.Range("B8").Value = "Nice"

$ws.Range("A9").Value = "Nice"
$ws.Range("B9").Value = "Nice"

$ws.Range("A10").Value = "Test"
$ws.Range("B10").Value = "TEst"

$ws.Range("A11").Value = "nice"
$ws.Range("B11").Value = "nice "

$ws.Range("A12").Value = "nice"
$ws.Range("B12").Value = "nice "

$ws.Range("A13").Value = "Non"
$ws.Range("B13").Value = "NOn"

$ws.Range("A14").Value = "no"
$ws.Range("B14").Value = "no"

$ws.Range("A15").Value = "no"
$ws.Range("B15").Value = "no"

$ws.Range("A16").Value = "as"
$ws.Range("B16").Value = "as"
